$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-D($Row, $Val) {
    $cell = $ws.Range("D$Row")
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}
function Set-E($Row, $Val) {
    $ws.Range("E$Row").Value = $Val
}

Set-D 2  "29.241.25"
Set-E 2  "  -0.81%  "

Set-D 3  "1.863.63"
Set-E 3  "  -0.90%  "

Set-E 4  "  -0.01%  "

Set-D 5  "0.7127"
Set-E 5  "  -0.99%  "

Set-D 6  "240.77"
Set-E 6  "  -0.93%  "

Set-D 7  "1.001"
Set-E 7  "  +0.05%  "

Set-D 8  "0.3084"
Set-E 8  "  -1.46%  "

Set-D 9  "0.07722"
Set-E 9  "  -1.94%  "

Set-D 10 "24.91"
Set-E 10 "  -1.46%  "

Set-D 11 "0.08347"
Set-E 11 "  +1.09%  "

Set-D 12 "1.860.94"
Set-E 12 "  -1.35%  "

Set-D 13 "5.203"
Set-E 13 "  -1.85%  "

Set-D 14 "0.7134"
Set-E 14 "  -2.70%  "

Set-D 15 "91.28"
Set-E 15 "  +0.11%  "

Set-D 16 "29.251.58"
Set-E 16 "  -0.80%  "

Set-D 17 "5.942"
Set-E 17 "  -0.35%  "

Set-D 18 "242.69"
Set-E 18 "  -1.98%  "

Set-D 19 "0.000007840"
Set-E 19 "  -0.90%  "

Set-D 20 "2.124.79"
Set-E 20 "  -0.22%  "

Set-E 21 "  -1.36%  "

Set-D 22 "1.0000"
Set-E 22 "  +0.12%  "

Set-D 23 "7.855"
Set-E 23 "  -1.78%  "

Set-D 24 "1.000"
Set-E 24 "  -0.02%  "

Set-D 25 "0.1595"
Set-E 25 "  +0.67%  "

Set-D 26 "163.33"
Set-E 26 "  -0.51%  "

Set-D 27 "8.892"
Set-E 27 "  -2.01%  "

Set-D 28 "18.49"
Set-E 28 "  +0.76%  "

Set-D 29 "1.343"
Set-E 29 "  -1.43%  "

Set-D 30 "1.497"
Set-E 30 "  -0.20%  "

Set-D 31 "4.411"
Set-E 31 "  +0.52%  "

Set-D 32 "4.247"
Set-E 32 "  +2.16%  "

Set-D 33 "0.05146"
Set-E 33 "  -2.99%  "

Set-D 34 "0.8224"
Set-E 34 "  +13.40%  "

Set-D 35 "1.932"
Set-E 35 "  -0.33%  "

Set-E 36 "  -3.26%  "

Set-D 37 "2.683"

Set-D 38 "0.01853"
Set-E 38 "  -0.97%  "

Set-D 39 "2.696"
Set-E 39 "  -1.28%  "

Set-D 40 "1.172.56"
Set-E 40 "  -7.31%  "

Set-D 41 "6.201"
Set-E 41 "  +1.49%  "

Set-D 42 "0.8941"
Set-E 42 "  -1.76%  "

Set-D 43 "72.81"
Set-E 43 "  -1.82%  "

Set-D 44 "0.9998"
Set-E 44 "  -0.02%  "

Set-D 45 "102.11"
Set-E 45 "  -1.50%  "

Set-D 46 "2.020.53"
Set-E 46 "  -0.42%  "

Set-D 47 "0.5206"
Set-E 47 "  -2.18%  "

Set-D 48 "1.792"
Set-E 48 "  +0.86%  "

Set-E 49 "  +0.03%  "

Set-D 50 "9.284"
Set-E 50 "  +0.05%  "

# Row 51: coin changed from Aptos to Frax
$ws.Range("B51").Value = "Frax"
$ws.Range("C51").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-D 51 "1.000"
Set-E 51 "  -0.20%  "
